$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Usuarios")

# --- Fix/standardize the passwords for every existing user (rows 2-8) ---
# They used to be a mix of "123" / "1234" / "mynorxd"; unify them all to
# a single, stronger password.
$newPassword = "Mynor123!"
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 3).Value = $newPassword
}

# --- Register the new user "mynor3" in row 9 with the same password ---
$ws.Cells.Item(9, 1).Value = "mynor3"
$ws.Cells.Item(9, 2).Value = "barriosmynitor@gmail.com"
$ws.Cells.Item(9, 3).Value = $newPassword

# --- Stray formatting left on G6 (underlined, otherwise empty cell) ---
$ws.Range("G6").Font.Underline = $true
$ws.Range("G6").Value = "a"
$ws.Range("G6").ClearContents()

# --- Restore the view/selection state ---
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("G6").Select()
